$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-03 Monday", "2025-02-04 Tuesday"),
    @("938×9=", "153×7="),
    @("744×3=", "437×3="),
    @("230×4=", "766×4="),
    @("113×3=", "255×7="),
    @("732×3=", "403×3="),
    @("678×4=", "425×5="),
    @("432×2=", "903×4="),
    @("984×7=", "413×7="),
    @("288×2=", "955×7="),
    @("783×8=", "135×6="),
    @("739×4=", "879×7="),
    @("635×3=", "406×4="),
    @("281×4=", "101×2="),
    @("646×7=", "131×2="),
    @("825×4=", "296×9="),
    @("777×2=", "972×9="),
    @("235×2=", "642×5="),
    @("629×3=", "921×7="),
    @("918×7=", "289×3="),
    @("475×5=", "746×3="),
    @("376×9=", "682×5="),
    @("242×2=", "498×5="),
    @("630×8=", "534×9="),
    @("757×8=", "179×3="),
    @("297×2=", "706×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
